$wb = $excel.ActiveWorkbook

$data = $wb.Worksheets.Item("Data")

# The "Data" sheet is protected, so temporarily unlock only the two cells
# we need to rewrite (keeps every other cell's style untouched).
$data.Range("A5:A6").Locked = $false

$data.Range("A5").Value = "CDC froncière"
$data.Range("A6").Value = "CDC locative"

# Restore the original cell formatting (the Locked=$false toggle above
# allocates a brand-new style record; copy the still-original style from
# an untouched row back onto A5:A6 so the style index is unchanged).
$data.Range("A2").Copy()
$data.Range("A5:A6").PasteSpecial(-4122, $null, $null, $null)

# Move the sheet's selection to A7 (matches the saved selection in the
# target file) without leaving "Data" as the active tab.
$data.Range("A7").Select()
$wb.Worksheets.Item("Prêts").Activate()
